$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Values --------------------------------------------------------------
# Header row (B1=0, C1=1)
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1

# Index column (A2:A6) = 0..4
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4

# Track-name column (B2:B6)
$ws.Range("B2").Value = 'Starboy'
$ws.Range("B3").Value = 'One Dance'
$ws.Range("B4").Value = 'Too Good'
$ws.Range("B5").Value = 'Fix You'
$ws.Range("B6").Value = 'Too Good'

# Audio-features column (C2:C6)
$ws.Range("C2").Value = '[{''acousticness'': 0.169, ''instrumentalness'': 5.49e-06, ''speechiness'': 0.281, ''tempo'': 186.041, ''id'': ''2IY7eOUDjw2ArKYxKa2jXc'', ''track_href'': ''https://api.spotify.com/v1/tracks/2IY7eOUDjw2ArKYxKa2jXc'', ''time_signature'': 4, ''duration_ms'': 230467, ''key'': 7, ''valence'': 0.477, ''danceability'': 0.682, ''uri'': ''spotify:track:2IY7eOUDjw2ArKYxKa2jXc'', ''mode'': 1, ''energy'': 0.592, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/2IY7eOUDjw2ArKYxKa2jXc'', ''loudness'': -7.033, ''type'': ''audio_features'', ''liveness'': 0.136}]'
$ws.Range("C3").Value = '[{''acousticness'': 0.00902, ''instrumentalness'': 0.00246, ''speechiness'': 0.0522, ''tempo'': 103.981, ''id'': ''12VWzyPDBCc8fqeWCAfNwR'', ''track_href'': ''https://api.spotify.com/v1/tracks/12VWzyPDBCc8fqeWCAfNwR'', ''time_signature'': 4, ''duration_ms'': 173987, ''key'': 1, ''valence'': 0.378, ''danceability'': 0.785, ''uri'': ''spotify:track:12VWzyPDBCc8fqeWCAfNwR'', ''mode'': 1, ''energy'': 0.617, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/12VWzyPDBCc8fqeWCAfNwR'', ''loudness'': -5.871, ''type'': ''audio_features'', ''liveness'': 0.351}]'
$ws.Range("C4").Value = '[{''acousticness'': 0.0606, ''instrumentalness'': 7.05e-05, ''speechiness'': 0.118, ''tempo'': 117.984, ''id'': ''7fJtPlEZKxu6gvkfBFc5tW'', ''track_href'': ''https://api.spotify.com/v1/tracks/7fJtPlEZKxu6gvkfBFc5tW'', ''time_signature'': 4, ''duration_ms'': 263373, ''key'': 7, ''valence'': 0.391, ''danceability'': 0.804, ''uri'': ''spotify:track:7fJtPlEZKxu6gvkfBFc5tW'', ''mode'': 1, ''energy'': 0.65, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/7fJtPlEZKxu6gvkfBFc5tW'', ''loudness'': -7.79, ''type'': ''audio_features'', ''liveness'': 0.102}]'
$ws.Range("C5").Value = '[{''acousticness'': 0.163, ''instrumentalness'': 0.00195, ''speechiness'': 0.0338, ''tempo'': 138.265, ''id'': ''7LVHVU3tWfcxj5aiPFEW4Q'', ''track_href'': ''https://api.spotify.com/v1/tracks/7LVHVU3tWfcxj5aiPFEW4Q'', ''time_signature'': 4, ''duration_ms'': 295533, ''key'': 3, ''valence'': 0.122, ''danceability'': 0.209, ''uri'': ''spotify:track:7LVHVU3tWfcxj5aiPFEW4Q'', ''mode'': 1, ''energy'': 0.418, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/7LVHVU3tWfcxj5aiPFEW4Q'', ''loudness'': -8.74, ''type'': ''audio_features'', ''liveness'': 0.113}]'
$ws.Range("C6").Value = '[{''acousticness'': 0.0606, ''instrumentalness'': 7.05e-05, ''speechiness'': 0.118, ''tempo'': 117.984, ''id'': ''7fJtPlEZKxu6gvkfBFc5tW'', ''track_href'': ''https://api.spotify.com/v1/tracks/7fJtPlEZKxu6gvkfBFc5tW'', ''time_signature'': 4, ''duration_ms'': 263373, ''key'': 7, ''valence'': 0.391, ''danceability'': 0.804, ''uri'': ''spotify:track:7fJtPlEZKxu6gvkfBFc5tW'', ''mode'': 1, ''energy'': 0.65, ''analysis_url'': ''https://api.spotify.com/v1/audio-analysis/7fJtPlEZKxu6gvkfBFc5tW'', ''loudness'': -7.79, ''type'': ''audio_features'', ''liveness'': 0.102}]'

# -- Formatting ------------------------------------------------------------
# Build the header/index style (bold, centered, thin box border) on one cell,
# then propagate it via copy/paste-special so the engine only ever mints a
# single extra cellXf (matches the target stylesheet exactly).
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B1").VerticalAlignment = -4160    # xlTop
$ws.Range("B1").Borders.LineStyle = 1        # xlContinuous (-> thin)

$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)          # xlPasteFormats
$ws.Range("A2:A6").PasteSpecial(-4122)       # xlPasteFormats
$excel.CutCopyMode = $false
